$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42602.583020833335
$ws.Range("B5").Value = "Bag"
$ws.Range("C5").Value = 6613
$ws.Range("D5").Value = 10320
$ws.Range("E5").Value = 1283
$ws.Range("F5").Value = 147
$ws.Range("G5").Value = 71
$ws.Range("H5").Value = 66
$ws.Range("I5").Value = 31
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 66
$ws.Range("M5").Value = 33
